$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Himel -10 (row 25 = Himel's daily entry)
$ws.Range("N25").Value = -10

# Taher +250 (row 27 = Taher's daily entry)
$ws.Range("N27").Value = 250

# Forhad +250 (row 28 = Forhad's daily entry)
$ws.Range("N28").Value = 250

# Who did the bazar (shopping) for this day's column
$ws.Range("N42").Value = "Himel"

# Bazar -750 -> Bazar TK spend for the day increases by 750
$ws.Range("N43").Value = 750

# Rice TK spend for the day increases by 2500
$ws.Range("N44").Value = 2500

# Match the author's final cell selection
$ws.Range("N26").Select()
